# TC19-Manage API Keys and Sign Out
# Collapse the 4-column header/row table down to the single
# "input_KeyName" column (what used to be column D), dropping the
# other three generated columns (A-C) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the surviving header text ("input_KeyName", previously in D1)
# into A1, keeping its existing style.
$ws.Range("A1").Value = "input_KeyName"

# Row 2's only remaining cell (previously D2) was already blank;
# make sure A2 ends up blank too, but still present as a tracked cell.
$ws.Range("A2").ClearContents()
$ws.Range("A2").Font.Bold = $false

# Drop columns B:C:D (their data + their custom column widths) so only
# column A remains, shifting D's old width onto A.
$ws.Columns("B:D").Delete()
$ws.Columns.Item(1).ColumnWidth = 14.1666666666667
